$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 848.1402832499998
$schedule.Range("F2").Value = 14.02348351934524

# --- Sheet: Detailed ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B21").Value = -5.01
$detailed.Range("B22").Value = -6.73332
$detailed.Range("B23").Value = -6.52883
$detailed.Range("C23").Value = "historical"
$detailed.Range("B24").Value = -7.88121
$detailed.Range("C24").Value = "historical"
$detailed.Range("B25").Value = 0
$detailed.Range("C25").Value = "historical"
$detailed.Range("B26").Value = 0.7
$detailed.Range("C26").Value = "historical"
$detailed.Range("B27").Value = 0.02274
$detailed.Range("B28").Value = -1.16054
$detailed.Range("B29").Value = -5.58973
$detailed.Range("B30").Value = -5.48208
$detailed.Range("B31").Value = -2.54301
$detailed.Range("B32").Value = -6
$detailed.Range("B33").Value = -2.97897
$detailed.Range("B34").Value = 36.0601
$detailed.Range("B35").Value = 24.59223
$detailed.Range("B36").Value = 9.847939999999999
$detailed.Range("B37").Value = 36.0601
$detailed.Range("B38").Value = 47.41716
$detailed.Range("B39").Value = 64.02958
$detailed.Range("B40").Value = 64.82379
$detailed.Range("B41").Value = 64.8901
$detailed.Range("B42").Value = 63.94385
$detailed.Range("B43").Value = 60.4431
$detailed.Range("B44").Value = 59.81888
$detailed.Range("B45").Value = 58.2941
$detailed.Range("B46").Value = 57.3
$detailed.Range("B47").Value = 58.81822
$detailed.Range("B48").Value = 61.0907
